$wb = $excel.ActiveWorkbook

# The underlying "Spreadsheet UO" upgrade changes the electrolyzer cell
# voltage loss input on the Input sheet. All downstream formulas
# (Calculations!B26, Calculations!B27, Output!B6, Output!B12) recalculate
# automatically from this single input change.
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B12").Value = 1.8503099269127579

$excel.CalculateFullRebuild()
